# US001 fully completed and Scrum_Backlog updated
#
# - Volunteer name correction: "Hitesh" -> "Harsh" (Tasks 1 & 2 of US001, cells D3/D4,
#   which share one entry in the shared-strings table).
# - Task 1 (G3) and Task 2 (G4) statuses flip from "In Progress" to "Closed",
#   completing User Story US001. The cells also pick up the same green/white
#   "Closed" look already used by the other closed cells in that column (e.g. G2),
#   so we copy that formatting across after updating the values.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Backlog")

# Correct the volunteer's name for Task 1 and Task 2 under US001.
$ws.Range("D3").Value = "Harsh"
$ws.Range("D4").Value = "Harsh"

# Mark Task 1 and Task 2 as Closed - US001 is now fully completed.
$ws.Range("G3").Value = "Closed"
$ws.Range("G4").Value = "Closed"

# Apply the existing "Closed" cell formatting (green fill / white bold font)
# from G2 to the two cells we just closed out.
$ws.Range("G2").Copy()
$ws.Range("G3").PasteSpecial(-4122)
$ws.Range("G4").PasteSpecial(-4122)
$excel.CutCopyMode = 0

# Restore the last-saved selection/view state.
$ws.Range("I7").Select()
